$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Permitir que ela se conecte",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Permitir que a aplicação se conecte",
    2)
